# Update cryptocurrency price/volume data (commit: Updated cryptos list on Fri Mar 24 20:25:01 UTC 2023 with GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.854.76"
$ws.Range("E2").Value = "  -1.28%  "

$ws.Range("D3").Value = "1.763.79"
$ws.Range("E3").Value = "  -2.72%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'323.22"
$ws.Range("E5").Value = "  -1.84%  "

$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").Value = "'0.4272"
$ws.Range("E7").Value = "  -3.66%  "

$ws.Range("D8").Value = "'0.3628"
$ws.Range("E8").Value = "  -2.81%  "

$ws.Range("D9").Value = "'0.07579"
$ws.Range("E9").Value = "  -1.50%  "

$ws.Range("D10").Value = "'42.77"
$ws.Range("E10").Value = "  -4.51%  "

$ws.Range("D11").Value = "'1.097"
$ws.Range("E11").Value = "  -2.42%  "

$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("D13").Value = "'20.81"
$ws.Range("E13").Value = "  -5.25%  "

$ws.Range("E14").Value = "  -3.36%  "

$ws.Range("D15").Value = "'7.283"
$ws.Range("E15").Value = "  -3.71%  "

$ws.Range("D16").Value = "1.751.31"
$ws.Range("E16").Value = "  -3.74%  "

$ws.Range("D17").Value = "'91.59"
$ws.Range("E17").Value = "  -1.86%  "

$ws.Range("D18").Value = "'0.00001069"
$ws.Range("E18").Value = "  -1.21%  "

$ws.Range("D19").Value = "'0.06376"
$ws.Range("E19").Value = "  -2.03%  "

$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("D21").Value = "'17.10"
$ws.Range("E21").Value = "  -2.35%  "

$ws.Range("D22").Value = "'5.923"
$ws.Range("E22").Value = "  -5.10%  "

$ws.Range("D23").Value = "27.882.03"
$ws.Range("E23").Value = "  -1.40%  "

$ws.Range("D24").Value = "'11.26"
$ws.Range("E24").Value = "  -3.80%  "

$ws.Range("D25").Value = "'2.124"
$ws.Range("E25").Value = "  +3.36%  "

$ws.Range("D26").Value = "'160.38"
$ws.Range("E26").Value = "  +3.16%  "

$ws.Range("D27").Value = "'20.32"
$ws.Range("E27").Value = "  -1.36%  "

$ws.Range("D28").Value = "1.956.95"
$ws.Range("E28").Value = "  -3.29%  "

$ws.Range("D29").Value = "'2.147"
$ws.Range("E29").Value = "  -7.57%  "

$ws.Range("D30").Value = "'125.15"
$ws.Range("E30").Value = "  -1.78%  "

$ws.Range("E31").Value = "  -6.38%  "

$ws.Range("D32").Value = "'3.686"
$ws.Range("E32").Value = "  +0.23%  "

$ws.Range("D33").Value = "'5.591"
$ws.Range("E33").Value = "  -4.59%  "

$ws.Range("D34").Value = "'0.08936"
$ws.Range("E34").Value = "  -3.07%  "

$ws.Range("D35").Value = "'12.28"

$ws.Range("D36").Value = "'0.02310"
$ws.Range("E36").Value = "  -1.61%  "

$ws.Range("D37").Value = "'0.2114"
$ws.Range("E37").Value = "  -2.59%  "

$ws.Range("D38").Value = "'0.06038"
$ws.Range("E38").Value = "  -2.61%  "

$ws.Range("D39").Value = "'0.6370"
$ws.Range("E39").Value = "  -2.89%  "

$ws.Range("D40").Value = "'4.997"
$ws.Range("E40").Value = "  -3.48%  "

$ws.Range("D41").Value = "'1.182"
$ws.Range("E41").Value = "  -1.42%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'7.926"
$ws.Range("E42").Value = "  -2.11%  "

$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").Value = "'0.9998"
$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("D44").Value = "'1.400"
$ws.Range("E44").Value = "  +1.07%  "

$ws.Range("D45").Value = "'13.44"
$ws.Range("E45").Value = "  -2.86%  "

$ws.Range("D46").Value = "'0.5895"
$ws.Range("E46").Value = "  -2.93%  "

$ws.Range("D47").Value = "'3.699"
$ws.Range("E47").Value = "  -1.66%  "

$ws.Range("D48").Value = "'1.997"
$ws.Range("E48").Value = "  -1.83%  "

$ws.Range("D49").Value = "'123.29"
$ws.Range("E49").Value = "  -2.61%  "

$ws.Range("D50").Value = "'1.187"
$ws.Range("E50").Value = "  +3.21%  "

$ws.Range("D51").Value = "'0.06838"
$ws.Range("E51").Value = "  -2.07%  "

